# Update "想去人数" (F column) values on sheets "展览" and "全部类型"
# as published by the gh-pages output generation at commit 456a3b4.

$wb = $excel.ActiveWorkbook

# Map of sheet name -> row -> new F value
$updates = @{
    "展览" = @{
        2  = 6950
        3  = 19
        8  = 114
        13 = 441
        14 = 24
        15 = 1813
        16 = 39
        17 = 3566
        20 = 78
        21 = 14
        22 = 24
        23 = 2186
        24 = 1
        25 = 227
        30 = 15
        31 = 148
        32 = 66
        33 = 42
    }
    "全部类型" = @{
        2  = 6950
        3  = 19
        9  = 114
        14 = 441
        15 = 24
        16 = 1813
        17 = 39
        18 = 3566
        21 = 78
        22 = 14
        23 = 24
        24 = 2186
        25 = 1
        26 = 227
        31 = 15
        32 = 148
        33 = 66
        34 = 43
    }
}

foreach ($sheetName in $updates.Keys) {
    $ws = $wb.Worksheets.Item($sheetName)
    $rows = $updates[$sheetName]
    foreach ($row in $rows.Keys) {
        $ws.Range("F$row").Value = $rows[$row]
    }
}
